$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Add a new cell with the reference link (and make it an actual hyperlink)
$ws.Range("D2").Value = "http://fabienpn.wordpress.com/2013/08/16/qt-thread-multiple-methods-with-sources/"
$ws.Hyperlinks.Add($ws.Range("D2"), "http://fabienpn.wordpress.com/2013/08/16/qt-thread-multiple-methods-with-sources/")
$ws.Columns.Item(4).ColumnWidth = 81.28515625

# Fix typo in the "communication between serial thread and gui thread" note
$ws.Range("C2").Value = "Look into communication between the serial thread and gui thread. Possible options are portected buffers/variables  or  futures."

# Move the active selection to C2
$ws.Range("C2").Select() | Out-Null
